$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Table 1 (rows 1-8): add an F1 Score column in column S
# ---------------------------------------------------------------------------
$ws.Range("S2").Value = "F1 Score"

$ws.Range("S3").Formula = "=(2*Q3*R3)/(Q3+R3)"
$ws.Range("S4:S7").Formula = "=(2*Q4*R4)/(Q4+R4)"
$ws.Range("S8").Formula = "=(2*Q8*R8)/(Q8+R8)"

# ---------------------------------------------------------------------------
# Table 2 (rows 11-25): add F1 Score columns in S (precision/recall block),
# X (micro block) and AC (macro block)
# ---------------------------------------------------------------------------
$ws.Range("S11").Value = "F1 Score"
$ws.Range("AC11").Value = "F1-Score"

$ws.Range("S12:S25").Formula = "=(2*Q12*R12)/(Q12+R12)"
$ws.Range("X12:X21").Formula = "=(2*V12*W12)/(V12+W12)"
$ws.Range("AC12:AC21").Formula = "=(2*AA12*AB12)/(AA12+AB12)"

# ---------------------------------------------------------------------------
# Restore the view/selection state recorded in the saved workbook
# ---------------------------------------------------------------------------
[void]$ws.Range("W10").Select()
